# Update the "想去人数" (want-to-go count) figures in column F across the
# four sheets of the workbook, as published by the gh-pages data refresh.

$wb = $excel.ActiveWorkbook

# 展览 (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 67
$ws1.Range("F3").Value  = 365
$ws1.Range("F5").Value  = 1268
$ws1.Range("F7").Value  = 2427
$ws1.Range("F8").Value  = 857
$ws1.Range("F9").Value  = 18402
$ws1.Range("F10").Value = 47
$ws1.Range("F11").Value = 1803
$ws1.Range("F12").Value = 645
$ws1.Range("F13").Value = 593
$ws1.Range("F18").Value = 63
$ws1.Range("F19").Value = 312
$ws1.Range("F21").Value = 85
$ws1.Range("F22").Value = 12
$ws1.Range("F23").Value = 73

# 演出 (Performances)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F7").Value  = 107
$ws2.Range("F9").Value  = 106
$ws2.Range("F11").Value = 8

# 本地生活 (Local life)
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 5846
$ws3.Range("F3").Value = 545

# 全部类型 (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 67
$ws4.Range("F3").Value  = 5846
$ws4.Range("F4").Value  = 545
$ws4.Range("F6").Value  = 365
$ws4.Range("F10").Value = 1268
$ws4.Range("F15").Value = 2427
$ws4.Range("F16").Value = 857
$ws4.Range("F17").Value = 18402
$ws4.Range("F18").Value = 47
$ws4.Range("F19").Value = 107
$ws4.Range("F22").Value = 1803
$ws4.Range("F23").Value = 645
$ws4.Range("F24").Value = 106
$ws4.Range("F25").Value = 593
$ws4.Range("F31").Value = 63
$ws4.Range("F32").Value = 8
$ws4.Range("F34").Value = 312
$ws4.Range("F39").Value = 85
$ws4.Range("F43").Value = 12
$ws4.Range("F48").Value = 73
